$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 22, pushing existing rows 22.. down to 24..
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

# New row 22: Cilantro "Primera" entry for the week of 2022-08-31
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value = "Ñuble"
$ws.Range("D22").Value = 44804
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 100112040
$ws.Range("G22").Value = "Cilantro"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 240
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = 750
$ws.Range("N22").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O22").Value = "Provincia de Diguillín"
$ws.Range("P22").Value = 750
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = "Hortaliza"

# New row 23: Cilantro "Segunda" entry for the week of 2022-08-31
$ws.Range("A23").Value = 7
$ws.Range("B23").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C23").Value = "Ñuble"
$ws.Range("D23").Value = 44804
$ws.Range("E23").Value = 16
$ws.Range("F23").Value = 100112040
$ws.Range("G23").Value = "Cilantro"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Segunda"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 600
$ws.Range("M23").Value = 600
$ws.Range("N23").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O23").Value = "Provincia de Diguillín"
$ws.Range("P23").Value = 600
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = "Hortaliza"
